# Appends new spectrum-analyzer marker log blocks to the "Results" sheet
# (rows 18-46), matching the pattern already present in rows 2-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainTextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $ws.Range($Address).Value = $Text
}

function Set-ProtectedTextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    # Force text storage so numeric-looking strings (e.g. "8.9470000000E+07")
    # are not auto-coerced into real numbers by Excel's type inference.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

$blocks = @(
    @{ Row = 18; Timestamp = "2025-06-26 11:54:15"; Freq = "8.9470000000E+07";    Power = "-1.000000E+03" },
    @{ Row = 23; Timestamp = "2025-06-26 11:56:48"; Freq = "8.9470000000E+07";    Power = "-1.000000E+03" },
    @{ Row = 28; Timestamp = "2025-06-26 12:10:37"; Freq = "8.9470000000E+07";    Power = "-1.000000E+03" },
    @{ Row = 33; Timestamp = "2025-06-26 12:14:24"; Freq = "2.241600000000E+09";  Power = "-8.8010E+01" },
    @{ Row = 38; Timestamp = "2025-06-26 12:30:07"; Freq = "2.203 GHz";           Power = "-88.88 dBm" },
    @{ Row = 43; Timestamp = "2025-06-26 12:31:52"; Freq = "2.248 GHz";           Power = "-88.99 dBm" }
)

foreach ($b in $blocks) {
    $r0 = $b.Row

    Set-PlainTextCell ("A" + $r0) "Timestamp"
    Set-PlainTextCell ("B" + $r0) $b.Timestamp

    $r2 = $r0 + 2
    Set-PlainTextCell ("A" + $r2) "Marker Frequency (Hz)"
    Set-PlainTextCell ("B" + $r2) "Marker Power (dBm)"

    $r3 = $r0 + 3
    Set-ProtectedTextCell ("A" + $r3) $b.Freq
    Set-ProtectedTextCell ("B" + $r3) $b.Power
}
